$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Completed intervention simulation for antimicrobial replacement
# intervention (ILRI probiotic trial data analysis) - updated scenario
# analysis results for Human.Capital.Approach (B) and Friction.Cost.Approach (C)
$ws.Range("B2").Value = -4385203656.2914505
$ws.Range("C2").Value = -4386228768.944093

$ws.Range("B3").Value = -4383336448.255414
$ws.Range("C3").Value = -4385386676.249562

$ws.Range("B4").Value = -4379602019.258911
$ws.Range("C4").Value = -4383702486.0027

$ws.Range("B5").Value = -4375120681.715998
$ws.Range("C5").Value = -4381681449.156696
